$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    "Sat Feb 17 22:37:38 EST 2024",
    "Sat Feb 17 22:37:48 EST 2024",
    "Sat Feb 17 22:37:57 EST 2024",
    "Sat Feb 17 22:38:06 EST 2024",
    "Sat Feb 17 22:38:15 EST 2024",
    "Sat Feb 17 22:38:24 EST 2024",
    "Sat Feb 17 22:38:35 EST 2024",
    "Sat Feb 17 22:38:44 EST 2024",
    "Sat Feb 17 22:38:54 EST 2024",
    "Sat Feb 17 22:39:03 EST 2024",
    "Sat Feb 17 22:39:13 EST 2024",
    "Sat Feb 17 22:39:22 EST 2024",
    "Sat Feb 17 22:39:31 EST 2024",
    "Sat Feb 17 22:39:40 EST 2024",
    "Sat Feb 17 22:39:49 EST 2024",
    "Sat Feb 17 22:39:58 EST 2024",
    "Sat Feb 17 22:40:07 EST 2024",
    "Sat Feb 17 22:40:16 EST 2024",
    "Sat Feb 17 22:40:26 EST 2024",
    "Sat Feb 17 22:40:35 EST 2024",
    "Sat Feb 17 22:40:44 EST 2024",
    "Sat Feb 17 22:40:53 EST 2024",
    "Sat Feb 17 22:41:02 EST 2024",
    "Sat Feb 17 22:41:11 EST 2024",
    "Sat Feb 17 22:41:20 EST 2024",
    "Sat Feb 17 22:41:30 EST 2024",
    "Sat Feb 17 22:41:39 EST 2024",
    "Sat Feb 17 22:41:48 EST 2024",
    "Sat Feb 17 22:41:57 EST 2024",
    "Sat Feb 17 22:42:06 EST 2024",
    "Sat Feb 17 22:42:17 EST 2024",
    "Sat Feb 17 22:42:26 EST 2024",
    "Sat Feb 17 22:42:35 EST 2024",
    "Sat Feb 17 22:42:44 EST 2024",
    "Sat Feb 17 22:42:53 EST 2024",
    "Sat Feb 17 22:43:02 EST 2024",
    "Sat Feb 17 22:43:12 EST 2024",
    "Sat Feb 17 22:43:22 EST 2024",
    "Sat Feb 17 22:43:31 EST 2024",
    "Sat Feb 17 22:43:40 EST 2024",
    "Sat Feb 17 22:43:50 EST 2024",
    "Sat Feb 17 22:43:59 EST 2024",
    "Sat Feb 17 22:44:08 EST 2024",
    "Sat Feb 17 22:44:17 EST 2024",
    "Sat Feb 17 22:44:26 EST 2024",
    "Sat Feb 17 22:44:35 EST 2024"
)

$startRow = 2
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("B$row").Value = $newDates[$i]
}
